$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.148.34"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").Value = "3.208.46"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'597.12"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.71%  "
$ws.Range("D6").Value = "'154.56"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").Value = "3.208.34"
$ws.Range("E8").Value = "  +1.95%  "
$ws.Range("D9").Value = "'0.535"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.72%  "
$ws.Range("D10").Value = "'0.161"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("D11").Value = "'6.10"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").Value = "'0.515"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.32%  "
$ws.Range("D13").Value = "'0.0000272"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.81%  "
$ws.Range("D14").Value = "'39.26"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.67%  "
$ws.Range("D15").Value = "3.729.67"
$ws.Range("E15").Value = "  +1.66%  "
$ws.Range("D16").Value = "66.050.97"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("D17").Value = "'7.44"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.68%  "
$ws.Range("D18").Value = "3.207.33"
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").Value = "'511.54"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").Value = "'15.35"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.73%  "
$ws.Range("D22").Value = "'0.746"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.05%  "
$ws.Range("D23").Value = "'15.30"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").Value = "'8.03"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.72%  "
$ws.Range("D25").Value = "'85.10"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "'9.42"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +5.91%  "
$ws.Range("D28").Value = "'3.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.26%  "
$ws.Range("D29").Value = "'2.28"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.09%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'6.91"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +10.25%  "
$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D31").Value = "'2.88"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.05%  "
$ws.Range("D32").Value = "'28.45"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.86%  "
$ws.Range("E33").Value = "  +3.42%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").Value = "'6.58"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.70%  "
$ws.Range("D36").Value = "'55.11"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("D37").Value = "'0.0910"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("D38").Value = "'485.98"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.82%  "
$ws.Range("D39").Value = "'0.0420"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("D40").Value = "'2.94"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("D41").Value = "'8.86"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.47%  "
$ws.Range("D42").Value = "'0.301"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.33%  "
$ws.Range("E43").Value = "  +3.36%  "
$ws.Range("D44").Value = "0.0₃0651"
$ws.Range("E44").Value = "  +10.10%  "
$ws.Range("D45").Value = "2.947.43"
$ws.Range("E45").Value = "  -3.48%  "
$ws.Range("D46").Value = "'2.44"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("D47").Value = "'28.55"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").Value = "'0.117"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("D50").Value = "'2.30"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("D51").Value = "'2.55"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.27%  "
